# SK렌터카 IFRS financial-summary refresh (rows 2-9 of company_list).
# Rows 2-6: figures rebased (likely units changed, e.g. KRW -> 100M KRW) and
#           several ratio-looking columns (AG/AH/AI) zeroed out.
# Rows 7-9: forecast years ("(E)") wiped back to only the label columns plus
#           AG/AH (=0); every other data cell is cleared entirely, matching
#           the source diff where those <c> elements disappear outright.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Cells($row, $values) {
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}

function Clear-Cells($row, $cols) {
    foreach ($col in $cols) {
        $ws.Range("$col$row").ClearContents()
    }
}

# ---- Row 2 ----------------------------------------------------------------
Set-Cells 2 @{
    D="5700"; E="437"; F="437"; G="245"; H="190"; I="190";
    K="8872"; L="6784"; M="2087"; N="2087";
    P="111"; Q="-705"; R="-89"; S="805"; T="126"; U="-832"; V="5800";
    W="7.67"; X="3.34"; Y="9.54"; Z="2.28"; AA="325.02"; AB="2054.73";
    AC="860"; AD="18.96"; AE="9426"; AF="1.73";
    AG="0"; AH="0"; AI="0"; AJ="22146300"
}
Clear-Cells 2 @("J","O")

# ---- Row 3 ----------------------------------------------------------------
Set-Cells 3 @{
    D="6321"; E="403"; F="403"; G="226"; H="171"; I="175"; J="-4";
    K="10338"; L="8063"; M="2274"; N="2253"; O="21";
    P="111"; Q="-532"; R="-400"; S="1145"; T="350"; U="-882"; V="6900";
    W="6.37"; X="2.71"; Y="8.06"; Z="1.79"; AA="354.57"; AB="2217.26";
    AC="790"; AD="15.06"; AE="10211"; AF="1.17";
    AG="0"; AH="0"; AI="0"; AJ="22146300"
}

# ---- Row 4 ----------------------------------------------------------------
Set-Cells 4 @{
    D="6476"; E="327"; F="327"; G="120"; H="68"; I="80"; J="-11";
    K="11576"; L="9233"; M="2344"; N="2334"; O="10";
    P="111"; Q="-720"; R="-138"; S="1085"; T="66"; U="-786"; V="8003";
    W="5.05"; X="1.05"; Y="3.48"; Z="0.62"; AA="393.98"; AB="2289.51";
    AC="360"; AD="23.18"; AE="10578"; AF="0.79";
    AG="0"; AH="0"; AI="0"; AJ="22146300"
}

# ---- Row 5 ----------------------------------------------------------------
Set-Cells 5 @{
    D="6368"; E="380"; F="380"; G="141"; H="94"; I="97"; J="-3";
    K="12122"; L="9702"; M="2419"; N="2413"; O="6";
    P="111"; Q="-496"; R="-249"; S="472"; T="197"; U="-693"; V="8476";
    W="5.97"; X="1.48"; Y="4.1"; Z="0.79"; AA="401.06"; AB="2364.05";
    AC="440"; AD="22.54"; AE="10937"; AF="0.91";
    AG="0"; AH="0"; AI="0"; AJ="22146300"
}

# ---- Row 6 ----------------------------------------------------------------
Set-Cells 6 @{
    D="6525"; E="218"; F="218"; G="42"; H="33"; I="38";
    K="11732"; L="9262"; M="2470"; N="2469";
    P="111"; Q="-253"; R="349"; S="-100"; T="65"; U="-319"; V="8038";
    W="3.35"; X="0.51"; Y="1.55"; Z="0.28"; AA="375.06"; AB="2407.62";
    AC="171"; AD="71.86"; AE="11193"; AF="1.1";
    AG="0"; AH="0"; AI="0"; AJ="22146300"
}

# ---- Rows 7-9: forecast years collapse to just AG/AH = 0 ------------------
$forecastClearCols = @("D","E","G","H","I","K","L","M","N","P","Q","R","S","T","U","W","X","Y","Z","AA","AC","AD","AE","AF","AI")
foreach ($row in 7..9) {
    Clear-Cells $row $forecastClearCols
    Set-Cells $row @{ AG="0"; AH="0" }
}
